# ------------------------------------------------------------------
# Adds the "2022-Q3" quarterly holdings sheet (inserted right after
# "总计") and updates the "总计" (summary) sheet with the new row.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" sheet: insert a new data row right under the header and
#    fill it with the 2022-Q3 summary figures. Everything below
#    shifts down by one row (handled by Rows.Insert()).
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# Re-apply the data-row formatting (border-free, centred index column)
# to the freshly inserted row by copying it from the row below, then
# overwrite the values.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 8
$summary.Range("D2").Value = 1.74

# Renumber the leading index column (0,1,2,...) for every data row so
# it keeps matching its row position after the insert.
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, positioned right after "总计" (i.e.
#    before the former second sheet, "2022-Q2"). Its layout mirrors
#    the other quarterly sheets, so duplicate "2022-Q2" (this keeps
#    all formatting -- fonts, borders, sheetPr, column styles --
#    intact, unlike building a blank sheet from scratch) and then
#    overwrite its contents with the 2022-Q3 numbers.
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $summary)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template ("2022-Q2") has 10 data rows (header + 9 funds); the
# 2022-Q3 table only has 9 (header + 8 funds), so drop the extra row.
$q3.Rows.Item(10).Delete()

# ---- header row (B1:H1) -----------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 2; $c -le 8; $c++) {
    $q3.Cells.Item(1, $c).Value = $headers[$c - 2]
}

# ---- data rows -----------------------------------------------------
# Columns: A index(n) B code(text) C name(text) D size(text)
#          E position(text) F ratio(text) G value(text) H rank(n)
$rows = @(
    @(0, "512980", "广发中证传媒ETF",              "44.76", "99.29", "2.97", "1.3294", 9),
    @(1, "160629", "鹏华中证传媒指数（LOF）A",       "6.41",  "94.58", "2.80", "0.1795", 9),
    @(2, "001628", "招商体育文化休闲股票A",          "2.23",  "92.42", "4.84", "0.1079", 10),
    @(3, "159805", "鹏华中证传媒ETF",                "1.71",  "98.37", "2.94", "0.0503", 9),
    @(4, "164818", "工银瑞信中证传媒指数（LOF）A",    "1.65",  "93.46", "2.77", "0.0457", 9),
    @(5, "015395", "招商体育文化休闲股票C",          "0.25",  "92.42", "4.84", "0.0121", 10),
    @(6, "010677", "工银瑞信中证传媒指数（LOF）C",    "0.21",  "93.46", "2.77", "0.0058", 9),
    @(7, "015675", "鹏华中证传媒指数（LOF）C",        "0.17",  "94.58", "2.80", "0.0048", 9)
)

# Text-typed numeric-looking columns need an explicit "@" (Text)
# number format before the assignment, otherwise the host silently
# parses them as numbers (and e.g. "001628" would lose its leading
# zero). Plain names in column C never look numeric, so they are
# left alone. After writing the value, the number format is put
# back so the cell keeps the same (style-less) look as its
# neighbours.
$textCols = @(2, 4, 5, 6, 7)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q3.Cells.Item($r, 1).Value = $data[0]

    foreach ($c in $textCols) {
        $q3.Cells.Item($r, $c).NumberFormat = "@"
    }

    $q3.Cells.Item($r, 2).Value = $data[1]
    $q3.Cells.Item($r, 3).Value = $data[2]
    $q3.Cells.Item($r, 4).Value = $data[3]
    $q3.Cells.Item($r, 5).Value = $data[4]
    $q3.Cells.Item($r, 6).Value = $data[5]
    $q3.Cells.Item($r, 7).Value = $data[6]

    foreach ($c in $textCols) {
        $q3.Cells.Item($r, $c).NumberFormat = "General"
        $q3.Cells.Item($r, $c).Style = "Normal"
    }

    $q3.Cells.Item($r, 8).Value = $data[7]
}

Write-Output "2022-Q3 sheet added and 总计 sheet updated"
